# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.898.88"
$ws.Range("E2").Value = "  +4.19%  "
$ws.Range("D3").Value = "2.779.68"
$ws.Range("E3").Value = "  +4.55%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'343.34"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "'115.36"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  +4.25%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("D10").Value = "'42.60"
$ws.Range("E10").Value = "  +6.94%  "
$ws.Range("D11").Value = "'0.0854"
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("D12").Value = "'20.01"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "'7.64"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "3.215.34"
$ws.Range("E15").Value = "  +4.64%  "
$ws.Range("D16").Value = "2.782.70"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "51.872.55"
$ws.Range("E18").Value = "  +4.26%  "
$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = "  +10.51%  "
$ws.Range("D20").Value = "'7.03"
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "'270.69"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'70.02"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +7.52%  "
$ws.Range("D26").Value = "'26.52"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'0.140"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "'34.64"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "'50.11"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "'18.98"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'2.10"
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").Value = "'4.95"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "'3.22"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").Value = "'0.0390"
$ws.Range("E40").Value = "  +12.23%  "
$ws.Range("D41").Value = "'2.65"
$ws.Range("E41").Value = "  +24.41%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'127.08"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'23.43"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.116"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").Value = "'2.32"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "2.066.31"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("D50").Value = "'0.899"
$ws.Range("E50").Value = "  +13.26%  "
$ws.Range("D51").Value = "'8.87"
$ws.Range("E51").Value = "  -1.53%  "
